$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.802.66'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '573.54'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.26%  '
$ws.Range('E6').Value = '  +0.48%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.537'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.33%  '
$ws.Range('D9').Value = '2.459.26'
$ws.Range('E9').Value = '  +0.60%  '
$ws.Range('E10').Value = '  +0.16%  '
$ws.Range('E11').Value = '  +0.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.29'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.14%  '
$ws.Range('E13').Value = '  +1.23%  '
$ws.Range('E15').Value = '  -0.96%  '
$ws.Range('D17').Value = '62.725.64'
$ws.Range('E17').Value = '  +0.16%  '
$ws.Range('D18').Value = '2.462.03'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.98'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.97'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '326.40'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.09%  '
$ws.Range('E22').Value = '  +10.26%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').Value = '  -0.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.13'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +19.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '65.62'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '656.74'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.65%  '
$ws.Range('D28').Value = '0.0₃0980'
$ws.Range('E28').Value = '  -0.70%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.998'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -15.39%  '
$ws.Range('E31').Value = '  -0.78%  '
$ws.Range('E32').Value = '  -2.55%  '
$ws.Range('E33').Value = '  -1.37%  '
$ws.Range('E34').Value = '  -2.77%  '
$ws.Range('E35').Value = '  -0.08%  '
$ws.Range('E36').Value = '  +2.32%  '
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('E38').Value = '  -1.49%  '
$ws.Range('E39').Value = '  -1.03%  '
$ws.Range('E40').Value = '  -2.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.69'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.37%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.79'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.00%  '
$ws.Range('E43').Value = '  -1.63%  '
$ws.Range('D44').Value = '0.0₆0310'
$ws.Range('E44').Value = '  -60.49%  '
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '153.34'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.54%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '15.22'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.57'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '20.45'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.87%  '
$ws.Range('E50').Value = '  +0.25%  '
$ws.Range('E51').Value = '  -1.28%  '
